$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 326
$ws.Range("I6").Value = 227.33333
$ws.Range("K6").Value = 681.99999
$ws.Range("M6").Value = -569.99999
$ws.Range("H9").Value = 100.25
$ws.Range("I9").Value = 83.666664
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 83.666664
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 85.333336
$ws.Range("N9").Value = -488
$ws.Range("H39").Value = 680.2759
$ws.Range("I39").Value = 339.35715
$ws.Range("J39").Value = 998.4666999999999
$ws.Range("K39").Value = 1018.07145
$ws.Range("L39").Value = 2995.4001
$ws.Range("M39").Value = -722.0714499999999
$ws.Range("N39").Value = -3587.4001
$ws.Range("H41").Value = 465.26666
$ws.Range("I41").Value = 292.42856
$ws.Range("K41").Value = 292.42856
$ws.Range("M41").Value = 147.57144
$ws.Range("H80").Value = 976.8182
$ws.Range("I80").Value = 1043.875
$ws.Range("J80").Value = 798
$ws.Range("K80").Value = 3131.625
$ws.Range("L80").Value = 2394
$ws.Range("M80").Value = -2133.625
$ws.Range("N80").Value = -4390
$ws.Range("H83").Value = 976.8182
$ws.Range("I83").Value = 1043.875
$ws.Range("J83").Value = 798
$ws.Range("K83").Value = 9394.875
$ws.Range("L83").Value = 7182
$ws.Range("M83").Value = -4402.875
$ws.Range("N83").Value = -17166
$ws.Range("H88").Value = 958.1818
$ws.Range("I88").Value = 222.83333
$ws.Range("J88").Value = 1840.6
$ws.Range("K88").Value = 222.83333
$ws.Range("L88").Value = 1840.6
$ws.Range("M88").Value = 183.16667
$ws.Range("N88").Value = -2652.6
$ws.Range("H91").Value = 958.1818
$ws.Range("I91").Value = 222.83333
$ws.Range("J91").Value = 1840.6
$ws.Range("K91").Value = 222.83333
$ws.Range("L91").Value = 1840.6
$ws.Range("M91").Value = 1181.16667
$ws.Range("N91").Value = -4648.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 4750
$ws.Range("I57").Value = 4750
$ws.Range("K57").Value = 4750
$ws.Range("M57").Value = -4266
$ws.Range("H80").Value = 39166.668
$ws.Range("J80").Value = 39166.668
$ws.Range("L80").Value = 39166.668
$ws.Range("N80").Value = -41162.668
$ws.Range("H83").Value = 39166.668
$ws.Range("J83").Value = 39166.668
$ws.Range("L83").Value = 117500.004
$ws.Range("N83").Value = -127484.004
$ws.Range("H88").Value = 1993.174
$ws.Range("J88").Value = 2618.1875
$ws.Range("L88").Value = 2618.1875
$ws.Range("N88").Value = -3430.1875
$ws.Range("H91").Value = 1993.174
$ws.Range("J91").Value = 2618.1875
$ws.Range("L91").Value = 2618.1875
$ws.Range("N91").Value = -5426.1875
$ws.Range("H110").Value = 4027.4614
$ws.Range("I110").Value = 786
$ws.Range("K110").Value = 786
$ws.Range("M110").Value = 1259
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14999.8
$ws.Range("J35").Value = 14999.8
$ws.Range("L35").Value = 14999.8
$ws.Range("N35").Value = -15619.8
$ws.Range("H86").Value = 2365.7222
$ws.Range("J86").Value = 4730.6
$ws.Range("L86").Value = 4730.6
$ws.Range("N86").Value = -6976.6
$ws.Range("H89").Value = 2365.7222
$ws.Range("J89").Value = 4730.6
$ws.Range("L89").Value = 23653
$ws.Range("N89").Value = -34885

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 237.125
$ws.Range("I10").Value = 128.14285
$ws.Range("K10").Value = 128.14285
$ws.Range("M10").Value = 10.85714999999999
$ws.Range("H19").Value = 564.9231
$ws.Range("I19").Value = 341.7143
$ws.Range("J19").Value = 1502.4
$ws.Range("K19").Value = 341.7143
$ws.Range("L19").Value = 1502.4
$ws.Range("M19").Value = -171.7143
$ws.Range("N19").Value = -1842.4
$ws.Range("H22").Value = 185.8
$ws.Range("H24").Value = 564.9231
$ws.Range("I24").Value = 341.7143
$ws.Range("J24").Value = 1502.4
$ws.Range("K24").Value = 341.7143
$ws.Range("L24").Value = 1502.4
$ws.Range("M24").Value = -171.7143
$ws.Range("N24").Value = -1842.4
$ws.Range("H31").Value = 3442
$ws.Range("I31").Value = 2487
$ws.Range("K31").Value = 2487
$ws.Range("M31").Value = -2192
$ws.Range("H34").Value = 3442
$ws.Range("I34").Value = 2487
$ws.Range("K34").Value = 2487
$ws.Range("M34").Value = -2285
$ws.Range("H43").Value = 10870
$ws.Range("J43").Value = 10870
$ws.Range("L43").Value = 10870
$ws.Range("N43").Value = -11238
$ws.Range("H62").Value = 3801.6667
$ws.Range("I62").Value = 3702.5
$ws.Range("K62").Value = 3702.5
$ws.Range("M62").Value = -3078.5
$ws.Range("H65").Value = 3801.6667
$ws.Range("I65").Value = 3702.5
$ws.Range("K65").Value = 18512.5
$ws.Range("M65").Value = -15392.5
$ws.Range("H101").Value = 10870
$ws.Range("J101").Value = 10870
$ws.Range("L101").Value = 10870
$ws.Range("N101").Value = -17360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 1100
$ws.Range("K69").Value = 3300
$ws.Range("M69").Value = -2489
$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 1100
$ws.Range("K72").Value = 9900
$ws.Range("M72").Value = -5844
$ws.Range("H80").Value = 4849
$ws.Range("J80").Value = 4800
$ws.Range("L80").Value = 14400
$ws.Range("N80").Value = -16272
$ws.Range("H83").Value = 4849
$ws.Range("J83").Value = 4800
$ws.Range("L83").Value = 43200
$ws.Range("N83").Value = -52560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 27296.166
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 36444.25
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 36444.25
$ws.Range("M43").Value = -8849
$ws.Range("N43").Value = -36746.25
$ws.Range("H97").Value = 1900.375
$ws.Range("I97").Value = 1900.375
$ws.Range("K97").Value = 1900.375
$ws.Range("M97").Value = -1404.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 450
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -227
$ws.Range("N55").Value = -846
$ws.Range("H82").Value = 1514.0834
$ws.Range("I82").Value = 854
$ws.Range("K82").Value = 854
$ws.Range("M82").Value = -493
$ws.Range("H85").Value = 1514.0834
$ws.Range("I85").Value = 854
$ws.Range("K85").Value = 854
$ws.Range("M85").Value = 394

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 41999.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 41999.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 41999.5
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -42553.5
$ws.Range("H62").Value = 2557.8
$ws.Range("I62").Value = 2597.25
$ws.Range("K62").Value = 2597.25
$ws.Range("M62").Value = -1973.25
$ws.Range("H65").Value = 2557.8
$ws.Range("I65").Value = 2597.25
$ws.Range("K65").Value = 12986.25
$ws.Range("M65").Value = -9866.25
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("H100").Value = 637.125
$ws.Range("I100").Value = 487.75
$ws.Range("J100").Value = 786.5
$ws.Range("K100").Value = 975.5
$ws.Range("L100").Value = 1573
$ws.Range("M100").Value = -434.5
$ws.Range("N100").Value = -2655
